try {
  [System.IO.File]::WriteAllText("C:\test.txt", "hi")
  Write-Host "file ok"
} catch {
  Write-Host "ERR:" $_.Exception.Message
}
